$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.349.34"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.500.79"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.70"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.15"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.71"
$ws.Range("E9").Value = "  +6.89%  "
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  +2.50%  "
$ws.Range("D12").Value = "4.097.83"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "3.501.80"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "64.304.23"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.32"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.56"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.74"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.580"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("D23").Value = "3.640.74"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.31"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.74"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("D34").Value = "3.529.76"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.30"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.42"
$ws.Range("E37").Value = "  +3.72%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "164.53"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0786"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.32"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "2.420.93"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("E51").Value = "  -0.07%  "
